$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Session" to "Neurology"
$ws.Name = "Neurology"

# New QR-scanner log rows to append (Student ID, Subject, Log Date, Log Time, Type, User)
$newRows = @(
    @("201850", "Neurology", "20/12/2025", "13:21:52", "Scan",   "emp17.farah.a.youssef@gmail.com"),
    @("201694", "Neurology", "20/12/2025", "13:21:56", "Scan",   "emp17.farah.a.youssef@gmail.com"),
    @("201790", "Neurology", "20/12/2025", "13:21:58", "Scan",   "emp17.farah.a.youssef@gmail.com"),
    @("201488", "Neurology", "20/12/2025", "13:22:13", "Manual", "emp17.farah.a.youssef@gmail.com"),
    @("201676", "Neurology", "20/12/2025", "13:22:16", "Scan",   "emp17.farah.a.youssef@gmail.com"),
    @("201619", "Neurology", "20/12/2025", "13:22:19", "Scan",   "emp17.farah.a.youssef@gmail.com"),
    @("201667", "Neurology", "20/12/2025", "13:22:22", "Scan",   "emp17.farah.a.youssef@gmail.com"),
    @("201675", "Neurology", "20/12/2025", "13:22:24", "Scan",   "emp17.farah.a.youssef@gmail.com"),
    @("201681", "Neurology", "20/12/2025", "13:22:29", "Scan",   "emp17.farah.a.youssef@gmail.com"),
    @("201561", "Neurology", "27/12/2025", "10:23:07", "Scan",   "emp17.farah.a.youssef@gmail.com"),
    @("212033", "Neurology", "27/12/2025", "10:23:47", "Manual", "emp17.farah.a.youssef@gmail.com"),
    @("201987", "Neurology", "27/12/2025", "10:24:26", "Manual", "emp17.farah.a.youssef@gmail.com"),
    @("201498", "Neurology", "27/12/2025", "10:24:38", "Manual", "emp17.farah.a.youssef@gmail.com")
)

$startRow = 80
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]

    # Column A holds digit-only IDs that must stay TEXT (not be coerced to a
    # number) - force the "@" text format before assigning the value, then
    # restore the default "Normal" style so no stray number-format style is
    # left behind on the cell.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $values[0]
    $cellA.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
    $ws.Cells.Item($r, 4).Value = $values[3]
    $ws.Cells.Item($r, 5).Value = $values[4]
    $ws.Cells.Item($r, 6).Value = $values[5]
}
